$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.9292377710648146, 0.7150252303766779, 0.6414142262493641, 3.195498068995942)"
$ws.Range("C2").Value = "NIG(1.618113666572148, 1.2352698058160239, 3.229232644310634, 6.646864796746801)"
$ws.Range("D2").Value = "NIG(0.9547437517562087, 0.6566310081163165, 1.9128786818540209, 3.047967231170968)"
$ws.Range("E2").Value = "NIG(1.9151233079232313, 1.6247049305186612, 3.1078689571255174, 5.438917940198555)"
